# Iteration 3 plan update: shift several deadlines from 06/05-10/05 week to
# 06/05-11/05, swap the "Nhóm" (group) numbers for items 2 and 3, and update
# the section title's date range accordingly.

$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Deadline column (column 4) updates, by table row (row 1 = header) ---
$t.Cell(2, 4).Range.Text  = "09/05"   # 08/05 -> 09/05
$t.Cell(3, 4).Range.Text  = "11/05"   # 09/05 -> 11/05
$t.Cell(4, 4).Range.Text  = "09/05"   # 07/05 -> 09/05
$t.Cell(5, 4).Range.Text  = "11/05"   # 10/05 -> 11/05
$t.Cell(6, 4).Range.Text  = "09/05"   # 08/05 -> 09/05
$t.Cell(7, 4).Range.Text  = "11/05"   # 10/05 -> 11/05
$t.Cell(8, 4).Range.Text  = "11/05"   # 09/05 -> 11/05
$t.Cell(9, 4).Range.Text  = "08/05"   # 07/05 -> 08/05
$t.Cell(10, 4).Range.Text = "11/05"   # 10/05 -> 11/05

# --- "Nhóm" column (column 2) swap for items 2 and 3 ---
$t.Cell(4, 2).Range.Text = "18"       # 17 -> 18
$t.Cell(6, 2).Range.Text = "17"       # 18 -> 17

# --- Title paragraph date range ---
$last = $d.Paragraphs.Count
$titleRange = $d.Paragraphs($last).Range
[void]$titleRange.Find.Execute("10/05", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "11/05", 2)
